$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# --- Update the per-row "time_taken" timestamps on the "data" sheet (F2:F64) ---
$newTimestamps = @(
    "2021-10-05 14:21:21.864648",
    "2021-10-05 14:21:21.864656",
    "2021-10-05 14:21:21.864660",
    "2021-10-05 14:21:21.864662",
    "2021-10-05 14:21:21.864665",
    "2021-10-05 14:21:21.864668",
    "2021-10-05 14:21:21.864671",
    "2021-10-05 14:21:21.864674",
    "2021-10-05 14:21:21.864677",
    "2021-10-05 14:21:21.864680",
    "2021-10-05 14:21:21.864682",
    "2021-10-05 14:21:21.864685",
    "2021-10-05 14:21:21.864687",
    "2021-10-05 14:21:21.864690",
    "2021-10-05 14:21:21.864692",
    "2021-10-05 14:21:21.864695",
    "2021-10-05 14:21:21.864697",
    "2021-10-05 14:21:21.864700",
    "2021-10-05 14:21:21.864703",
    "2021-10-05 14:21:21.864705",
    "2021-10-05 14:21:21.864708",
    "2021-10-05 14:21:21.864710",
    "2021-10-05 14:21:21.864713",
    "2021-10-05 14:21:21.864716",
    "2021-10-05 14:21:21.864718",
    "2021-10-05 14:21:21.864721",
    "2021-10-05 14:21:21.864724",
    "2021-10-05 14:21:21.864726",
    "2021-10-05 14:21:21.864729",
    "2021-10-05 14:21:21.864732",
    "2021-10-05 14:21:21.864734",
    "2021-10-05 14:21:21.864737",
    "2021-10-05 14:21:21.864740",
    "2021-10-05 14:21:21.864742",
    "2021-10-05 14:21:21.864745",
    "2021-10-05 14:21:21.864747",
    "2021-10-05 14:21:21.864750",
    "2021-10-05 14:21:21.864753",
    "2021-10-05 14:21:21.864755",
    "2021-10-05 14:21:21.864758",
    "2021-10-05 14:21:21.864761",
    "2021-10-05 14:21:21.864763",
    "2021-10-05 14:21:21.864766",
    "2021-10-05 14:21:21.864768",
    "2021-10-05 14:21:21.864771",
    "2021-10-05 14:21:21.864773",
    "2021-10-05 14:21:21.864776",
    "2021-10-05 14:21:21.864778",
    "2021-10-05 14:21:21.864781",
    "2021-10-05 14:21:21.864784",
    "2021-10-05 14:21:21.864786",
    "2021-10-05 14:21:21.864789",
    "2021-10-05 14:21:21.864792",
    "2021-10-05 14:21:21.864795",
    "2021-10-05 14:21:21.864797",
    "2021-10-05 14:21:21.864800",
    "2021-10-05 14:21:21.864802",
    "2021-10-05 14:21:21.864805",
    "2021-10-05 14:21:21.864807",
    "2021-10-05 14:21:21.864810",
    "2021-10-05 14:21:21.864812",
    "2021-10-05 14:21:21.864815",
    "2021-10-05 14:21:21.864818"
)
for ($i = 0; $i -lt $newTimestamps.Count; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimestamps[$i]
}

# --- Add the new "metadata" sheet after "data" ---
$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = "metadata"

# Header row (bold / centered / bordered, matching the "data" sheet's header style)
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

$dataSheet.Range("B1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)

# Data row
$ws.Range("A2").Value = 0
$dataSheet.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Range("B2").Value = "Laterality disorders and isomerism"
$ws.Range("C2").Value = 549
$ws.Range("D2").Value = "'1.45"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "2021-05-04T10:22:47.518757Z"
$ws.Range("F2").Value = "2021-10-05 14:21:21.861056"
$ws.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/549/?format=json"

Write-Output "done"
